# Weekly driver report update for 2025-04-19
# Updates the "Bad Drivers" summary numbers and re-sorts / refreshes the
# "Good Drivers" table (rows 13-25) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: assign a literal text value to a cell without Excel's "smart"
# text-to-number/date autoconversion touching the cell, and without
# creating any new cell style (keeps the original style index intact).
# We do this by writing a text FORMULA that evaluates to the literal
# string, then converting the formula to a static value via
# PasteSpecial(xlPasteValues) - this mirrors exactly what Excel does
# internally and never perturbs number formats / styles.
# ---------------------------------------------------------------------
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# -----------------------------------------------------------------
# 1) "Bad Drivers" table - refreshed Critical Minutes / Roaming %
# -----------------------------------------------------------------
$ws.Range("C3").Value = 259
$ws.Range("D3").Value = 92.5
$ws.Range("C4").Value = 1355
$ws.Range("D4").Value = 96.90000000000001
$ws.Range("C5").Value = 1614

# -----------------------------------------------------------------
# 2) "Good Drivers" table - the 13 data rows (13-25) get re-sorted
#    (newest driver vintage moves toward the bottom) and several
#    rows pick up refreshed sample counts / percentages. Some rows
#    also gain or lose their "Driver Vintage" date.
#
#    Stage the current A:E block into the unused F:J columns (still
#    inside the sheet's existing A1:J30 used range) so every source
#    cell stays readable/unchanged while we overwrite A:E in place.
# -----------------------------------------------------------------
$ws.Range("A13:E25").Copy($ws.Range("F13:J25"))

# staged row -> (Adapter-Driver text, Driver Vintage text) live in
# columns F (driver) and J (vintage) of the staging block.
# NOTE: always materialise the staged cell into a variable before calling
# a method on it - chaining a method directly off a function's return
# value (e.g. `Get-StagedDriver($r).Copy(...)`) does not reliably resolve
# against the right object in this host, so every call site below goes
# through an intermediate `$srcCell` variable.
function Get-StagedDriver($row) { return $ws.Cells.Item($row, 6) }   # col F
function Get-StagedVintage($row) { return $ws.Cells.Item($row, 10) } # col J

# target row => source (staged) row that its Adapter-Driver name/vintage came from
$rowMap = @{
    13 = 18
    14 = 15
    15 = 16
    16 = 17
    17 = 14
    18 = 19
    19 = 24
    20 = 13
    21 = 20
    22 = 22
    23 = 21
    24 = 23
    25 = 25
}

# target row => new Total Samples (col B) / Good Roaming % (col D)
$bValues = @{ 13 = 96526;  14 = 328411; 15 = 143808; 16 = 287148; 17 = 69578
              18 = 67111;  19 = 13016;  20 = 18721;  21 = 66577;  22 = 26241
              23 = 15730;  24 = 88435;  25 = 46270 }
$dValues = @{ 13 = 99.90000000000001; 14 = 99.90000000000001; 15 = 99.90000000000001
              16 = 99.90000000000001; 17 = 99.90000000000001; 18 = 100
              19 = 100; 20 = 99.90000000000001; 21 = 100; 22 = 100
              23 = 99.90000000000001; 24 = 99.90000000000001; 25 = 100 }

# target row => Driver Vintage handling: "clear" (blank it out), "copy"
# (pull the date text from the given staged row's vintage cell), or a
# literal new date string that didn't exist anywhere in the sheet before.
$eClear = @(13, 14, 15, 16, 17, 18, 19)
$eCopyFrom = @{ 20 = 13; 22 = 19; 23 = 21; 24 = 23; 25 = 25 }
$eLiteral = @{ 21 = "2024-05-09" }

foreach ($row in 13..25) {
    $srcRow = $rowMap[$row]

    # Column A - Adapter-Driver name (copy verbatim from its staged source row)
    $srcCell = Get-StagedDriver($srcRow)
    $dstCell = $ws.Cells.Item($row, 1)
    $srcCell.Copy($dstCell)

    # Column B - Total Samples
    $ws.Cells.Item($row, 2).Value = $bValues[$row]

    # Column D - Good Roaming Calculation (%)
    $ws.Cells.Item($row, 4).Value = $dValues[$row]

    # Column E - Driver Vintage
    if ($eClear -contains $row) {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.Clear()
    } elseif ($eCopyFrom.ContainsKey($row)) {
        $vSrcRow = $eCopyFrom[$row]
        $vSrcCell = Get-StagedVintage($vSrcRow)
        $vDstCell = $ws.Cells.Item($row, 5)
        $vSrcCell.Copy($vDstCell)
    } elseif ($eLiteral.ContainsKey($row)) {
        $eDstCell = $ws.Cells.Item($row, 5)
        Set-TextValue $eDstCell $eLiteral[$row]
    }
}

# Remove the staging block now that every row has been rebuilt.
$ws.Range("F13:J25").Clear()

Write-Output "driver summary refreshed"
